$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 130800693
$ws.Range("Q9").Value = 489211
$ws.Range("R9").Value = 6720001
$ws.Range("Z9").Value = '10:42'
$ws.Range("AB9").Value = '10:42'

# Row 10
$ws.Range("A10").Value = 130800672
$ws.Range("Q10").Value = 489170
$ws.Range("R10").Value = 6720188
$ws.Range("Z10").Value = '10:05'
$ws.Range("AB10").Value = '10:05'

# Row 11
$ws.Range("A11").Value = 130800627
$ws.Range("B11").Value = 57881
$ws.Range("E11").Value = 100049
$ws.Range("F11").Value = 'Spillkråka'
$ws.Range("G11").Value = 'Dryocopus martius'
$ws.Range("H11").Value = '(Linnaeus, 1758)'
$ws.Range("M11").Value = 'färska spår'
$ws.Range("Q11").Value = 489163
$ws.Range("R11").Value = 6719987
$ws.Range("Z11").Value = '11:08'
$ws.Range("AB11").Value = '11:08'

# Row 12
$ws.Range("A12").Value = 130800622
$ws.Range("M12").Value = 'äldre spår'
$ws.Range("Q12").Value = 489171
$ws.Range("R12").Value = 6720048
$ws.Range("Z12").Value = '10:24'
$ws.Range("AB12").Value = '10:24'

# Row 13
$ws.Range("A13").Value = 130800706
$ws.Range("B13").Value = 79243
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = 'Garnlav'
$ws.Range("G13").Value = 'Alectoria sarmentosa'
$ws.Range("H13").Value = '(Ach.) Ach.'
$ws.Range("Q13").Value = 489143
$ws.Range("R13").Value = 6720009
$ws.Range("Z13").Value = '11:13'
$ws.Range("AB13").Value = '11:13'
$ws.Range("M13").ClearContents()

# Row 14
$ws.Range("A14").Value = 130800721
$ws.Range("Q14").Value = 489189
$ws.Range("R14").Value = 6720073
$ws.Range("Z14").Value = '11:34'
$ws.Range("AB14").Value = '11:34'

# Row 26
$ws.Range("A26").Value = 130800606
$ws.Range("B26").Value = 92179
$ws.Range("D26").Value = 'VU'
$ws.Range("E26").Value = 2062
$ws.Range("F26").Value = 'Ulltickeporing'
$ws.Range("G26").Value = 'Skeletocutis brevispora'
$ws.Range("H26").Value = 'Niemelä'
$ws.Range("Q26").Value = 489169
$ws.Range("R26").Value = 6720185
$ws.Range("Z26").Value = '10:07'
$ws.Range("AB26").Value = '10:07'

# Row 27
$ws.Range("A27").Value = 130800725
$ws.Range("B27").Value = 79243
$ws.Range("D27").Value = 'NT'
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = 'Garnlav'
$ws.Range("G27").Value = 'Alectoria sarmentosa'
$ws.Range("H27").Value = '(Ach.) Ach.'
$ws.Range("Q27").Value = 489332
$ws.Range("R27").Value = 6720052
$ws.Range("Z27").Value = '11:50'
$ws.Range("AB27").Value = '11:50'

# Row 34
$ws.Range("A34").Value = 130800712
$ws.Range("Q34").Value = 489143
$ws.Range("R34").Value = 6720114
$ws.Range("Z34").Value = '11:21'
$ws.Range("AB34").Value = '11:21'

# Row 35
$ws.Range("A35").Value = 130800669
$ws.Range("Q35").Value = 489170
$ws.Range("R35").Value = 6720167
$ws.Range("Z35").Value = '10:03'
$ws.Range("AB35").Value = '10:03'

# Row 54
$ws.Range("A54").Value = 130800680
$ws.Range("Q54").Value = 489131
$ws.Range("R54").Value = 6720110
$ws.Range("Z54").Value = '10:16'
$ws.Range("AB54").Value = '10:16'

# Row 55
$ws.Range("A55").Value = 130800637
$ws.Range("B55").Value = 8451
$ws.Range("D55").Value = 'LC'
$ws.Range("E55").Value = 106545
$ws.Range("F55").Value = 'Mindre märgborre'
$ws.Range("G55").Value = 'Tomicus minor'
$ws.Range("H55").Value = '(Hartig, 1834)'
$ws.Range("M55").Value = 'färska gnagspår'
$ws.Range("Q55").Value = 489091
$ws.Range("R55").Value = 6720083
$ws.Range("Z55").Value = '11:19'
$ws.Range("AB55").Value = '11:19'
$ws.Range("AF55").Value = ''

# Row 56
$ws.Range("A56").Value = 130800683
$ws.Range("B56").Value = 79243
$ws.Range("D56").Value = 'NT'
$ws.Range("E56").Value = 6425
$ws.Range("F56").Value = 'Garnlav'
$ws.Range("G56").Value = 'Alectoria sarmentosa'
$ws.Range("H56").Value = '(Ach.) Ach.'
$ws.Range("Q56").Value = 489169
$ws.Range("R56").Value = 6720075
$ws.Range("Z56").Value = '10:20'
$ws.Range("AB56").Value = '10:20'
$ws.Range("M56").ClearContents()
$ws.Range("AF56").ClearContents()
